# Commit: "fixing errors with UI"
# Adds a new "Sample Item" data row (row 4) to both the "Monthly Costs"
# sheet and the "Pricing" sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Monthly Costs" ----
$ws1 = $wb.Worksheets.Item("Monthly Costs")

# Column A holds text-formatted dates (e.g. "2024-10-10") in the existing
# rows, so force a text number format before assigning the value to avoid
# Excel auto-converting the string into a date serial number.
$ws1.Range("A4").NumberFormat = "@"
$ws1.Range("A4").Value = "2024-10-10"
$ws1.Range("C4").Value = "Sample Item"
$ws1.Range("D4").Value = 100
$ws1.Range("E4").Value = 12.5

# ---- Sheet 3: "Pricing" ----
$ws3 = $wb.Worksheets.Item("Pricing")

$ws3.Range("A4").Value = "Sample Item"
$ws3.Range("B4").Value = 12.5
